$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content updates for the UKIM (people with significant control) template:
# shorten / reword the column header text.
$ws.Range("B1").Value = "Residential address"
$ws.Range("D1").Value = "National Insurance number"
$ws.Range("E1").Value = "Identification number if no National Insurance number (eg passport number, driver's licence, national identity card)"

# Columns F:Z on row 1 previously carried a slightly different (but still
# bold) style than A1:E1. Re-apply A1's formatting across the rest of the
# header row so every cell shares a single consistent style.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:Z1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
